$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.929.54'
$ws.Range("E2").Value = '  -3.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.354.01'
$ws.Range("E3").Value = '  -2.93%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.36'
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.77'
$ws.Range("E6").Value = '  +1.45%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.94'
$ws.Range("E9").Value = '  +1.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.416'
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.936.64'
$ws.Range("E12").Value = '  -2.70%  '
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.24'
$ws.Range("E14").Value = '  -1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.375.84'
$ws.Range("E15").Value = '  -2.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.052.96'
$ws.Range("E17").Value = '  -3.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.35'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.19'
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.86'
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '373.51'
$ws.Range("E21").Value = '  -3.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.562'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.05'
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.531.70'
$ws.Range("E25").Value = '  -1.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000108'
$ws.Range("E26").Value = '  -5.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.176'
$ws.Range("E27").Value = '  -3.73%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.41'
$ws.Range("E29").Value = '  -3.53%  '
$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.70'
$ws.Range("E32").Value = '  -4.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.89'
$ws.Range("E33").Value = '  -1.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.28'
$ws.Range("E34").Value = '  -3.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.38'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '169.75'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -4.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.77'
$ws.Range("E38").Value = '  -3.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '29.84'
$ws.Range("E39").Value = '  -6.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.391.94'
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0755'
$ws.Range("E41").Value = '  -1.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.33'
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("E43").Value = '  -3.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.30'
$ws.Range("E44").Value = '  -1.55%  '
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.14'
$ws.Range("E45").Value = '  -3.91%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.61'
$ws.Range("E46").Value = '  -6.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.512.45'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.87'
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.70'
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0261'
$ws.Range("E51").Value = '  -2.50%  '
